$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.98"
$ws.Range("E2").Value = "'2.43%"
$ws.Range("G2").Value = "'13"
$ws.Range("D3").Value = "'42.88"
$ws.Range("E3").Value = "'5.61%"
$ws.Range("G3").Value = "'13"
$ws.Range("D4").Value = "'5.041"
$ws.Range("E4").Value = "'0.63%"
$ws.Range("G4").Value = "'13"
$ws.Range("D5").Value = "'0.07674"
$ws.Range("E5").Value = "'3.82%"
$ws.Range("G5").Value = "'13"
$ws.Range("E6").Value = "'3.26%"
$ws.Range("G6").Value = "'13"
$ws.Range("D7").Value = "'1.030"
$ws.Range("E7").Value = "'11.26%"
$ws.Range("G7").Value = "'13"
$ws.Range("E8").Value = "'2.20%"
$ws.Range("G8").Value = "'13"
$ws.Range("D9").Value = "'0.1219"
$ws.Range("E9").Value = "'1.77%"
$ws.Range("G9").Value = "'13"
$ws.Range("E10").Value = "'2.27%"
$ws.Range("G10").Value = "'13"
$ws.Range("D11").Value = "'0.09056"
$ws.Range("E11").Value = "'2.68%"
$ws.Range("G11").Value = "'13"
$ws.Range("D12").Value = "'0.04154"
$ws.Range("E12").Value = "'-5.49%"
$ws.Range("G12").Value = "'13"
$ws.Range("D13").Value = "'0.1044"
$ws.Range("E13").Value = "'-1.14%"
$ws.Range("G13").Value = "'13"
$ws.Range("D14").Value = "'0.001280"
$ws.Range("E14").Value = "'1.26%"
$ws.Range("G14").Value = "'13"
$ws.Range("D15").Value = "'0.005947"
$ws.Range("E15").Value = "'2.11%"
$ws.Range("G15").Value = "'13"
$ws.Range("E16").Value = "'1,895.60%"
$ws.Range("G16").Value = "'13"
$ws.Range("D17").Value = "'3.320"
$ws.Range("E17").Value = "'-1.07%"
$ws.Range("G17").Value = "'13"
$ws.Range("D18").Value = "'4.406"
$ws.Range("E18").Value = "'2.63%"
$ws.Range("G18").Value = "'13"
$ws.Range("D19").Value = "'0.3339"
$ws.Range("E19").Value = "'1.99%"
$ws.Range("G19").Value = "'13"
$ws.Range("D20").Value = "'8.433"
$ws.Range("E20").Value = "'7.23%"
$ws.Range("G20").Value = "'13"
$ws.Range("D21").Value = "'0.1374"
$ws.Range("E21").Value = "'-1.17%"
$ws.Range("G21").Value = "'13"
$ws.Range("D22").Value = "'0.2987"
$ws.Range("E22").Value = "'6.52%"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'0.04142"
$ws.Range("E23").Value = "'5.59%"
$ws.Range("G23").Value = "'13"
$ws.Range("D24").Value = "'0.001270"
$ws.Range("E24").Value = "'0.37%"
$ws.Range("G24").Value = "'13"
$ws.Range("D25").Value = "'0.004505"
$ws.Range("E25").Value = "'18.40%"
$ws.Range("G25").Value = "'13"
$ws.Range("D26").Value = "'0.0001347"
$ws.Range("E26").Value = "'9.49%"
$ws.Range("G26").Value = "'13"
$ws.Range("G27").Value = "'13"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("D38").Value = "'0.02451"
$ws.Range("E38").Value = "'4.81%"
$ws.Range("G38").Value = "'13"
$ws.Range("D39").Value = "'0.05277"
$ws.Range("E39").Value = "'3.52%"
$ws.Range("G39").Value = "'13"
$ws.Range("D40").Value = "'0.006077"
$ws.Range("E40").Value = "'-1.57%"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.007657"
$ws.Range("E41").Value = "'-2.21%"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.1348"
$ws.Range("E42").Value = "'4.07%"
$ws.Range("G42").Value = "'13"
$ws.Range("D43").Value = "'0.007353"
$ws.Range("E43").Value = "'-0.47%"
$ws.Range("G43").Value = "'13"
$ws.Range("D44").Value = "'0.007377"
$ws.Range("E44").Value = "'0.18%"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.3033"
$ws.Range("E45").Value = "'3.07%"
$ws.Range("G45").Value = "'13"
$ws.Range("D46").Value = "'0.00006573"
$ws.Range("E46").Value = "'7.55%"
$ws.Range("G46").Value = "'13"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.27%"
$ws.Range("G47").Value = "'13"
$ws.Range("D48").Value = "'0.04586"
$ws.Range("E48").Value = "'-1.61%"
$ws.Range("G48").Value = "'13"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("G49").Value = "'13"
$ws.Range("D50").Value = "'0.00002095"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("G50").Value = "'13"
$ws.Range("D51").Value = "'0.0001995"
$ws.Range("E51").Value = "'-0.27%"
$ws.Range("G51").Value = "'13"
